$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.888.89"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.106.92"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "389.49"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.24"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.24"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0867"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "3.593.08"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.88"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "3.106.39"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.986"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("E18").Value = "  -3.93%  "
$ws.Range("D19").Value = "52.053.94"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.49"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.68"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.24"
$ws.Range("E26").Value = "  +5.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.13"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.109"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.39"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.79"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.39"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.299"
$ws.Range("E39").Value = "  +9.38%  "
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.10"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.71"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.74"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.46"
$ws.Range("E48").Value = "  +3.96%  "
$ws.Range("D49").Value = "2.058.32"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "3.413.09"
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.208"
$ws.Range("E51").Value = "  +7.38%  "
